$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = '{''control,localization'': 1, ''canbus,drivers,localization'': 1, ''canbus,control'': 4, ''common,control'': 2, ''canbus,common,control'': 1, ''control,drivers'': 1, ''canbus,common,control,localization,planning'': 3, ''common,control,localization'': 3, ''canbus,common,control,localization'': 2, ''canbus,control,planning'': 1, ''common,drivers'': 3, ''control,planning'': 1, ''canbus,drivers'': 1}'
$ws.Range("H2").Value = 646

$ws.Range("G3").Value = '{''canbus,control,perception'': 1, ''drivers,perception'': 6, ''localization,planning,prediction'': 1, ''map,planning'': 3, ''common,perception'': 3, ''canbus,common,control,dreamview,localization,perception,planning,prediction,routing'': 2, ''canbus,common,localization,planning,prediction'': 1, ''common,localization,planning,prediction'': 1, ''common,control,dreamview,localization,perception,planning,prediction'': 2, ''canbus,common,control,dreamview,localization,perception,planning,prediction'': 5, ''drivers,map,perception,planning'': 1, ''map,prediction'': 1, ''canbus,control,map,perception,planning,prediction,routing'': 1, ''control,planning'': 3, ''drivers,planning'': 1, ''planning,prediction'': 6, ''control,planning,prediction'': 1, ''common,control,perception,planning'': 2, ''control,perception'': 3, ''control,localization'': 2, ''control,localization,planning'': 2, ''common,map'': 2, ''common,planning'': 1, ''common,localization'': 1, ''drivers,map,perception'': 1, ''canbus,drivers'': 2, ''map,perception'': 3, ''calibration,dreamview,planning'': 1, ''dreamview,map'': 1, ''map,perception,planning'': 1, ''drivers,perception,planning'': 2, ''dreamview,routing'': 1, ''canbus,control'': 1, ''canbus,perception'': 1}'
$ws.Range("H3").Value = 1345

$ws.Range("G4").Value = '{''control,localization,perception'': 2, ''common,planning'': 2, ''canbus,common,control,drivers,localization,perception,planning,prediction'': 1, ''common,control,localization,perception,planning,prediction'': 3, ''canbus,common,control,dreamview,localization,perception,planning,prediction'': 3, ''common,control,dreamview,localization,perception,planning,prediction'': 4, ''perception,planning,prediction'': 1, ''localization,map'': 2, ''drivers,localization,map,perception'': 3, ''perception,planning,prediction,routing'': 1, ''perception,planning'': 4, ''map,perception'': 1, ''planning,prediction'': 2, ''localization,planning,prediction'': 2, ''perception,prediction'': 2, ''control,perception,planning,prediction'': 1, ''calibration,dreamview,perception'': 1, ''common,control,perception,planning'': 1, ''common,drivers,planning'': 1, ''common,control,drivers,planning'': 1, ''dreamview,drivers,localization,map,perception'': 1, ''control,localization'': 4, ''common,localization,perception,planning,prediction'': 2, ''common,localization,planning,prediction'': 1, ''control,drivers'': 2, ''control,dreamview,drivers'': 1, ''localization,map,perception'': 2, ''drivers,localization,perception'': 1, ''dreamview,map,planning,prediction'': 2, ''common,control'': 1, ''control,planning'': 4, ''localization,planning'': 6, ''dreamview,map'': 1, ''drivers,routing'': 1, ''common,perception'': 2, ''drivers,map,perception'': 1, ''drivers,perception,planning'': 1, ''control,prediction'': 1, ''common,map'': 1, ''common,map,perception'': 1, ''common,map,perception,planning'': 1, ''drivers,map'': 1, ''canbus,control'': 1, ''drivers,perception'': 3, ''canbus,perception'': 1, ''map,planning'': 1, ''localization,perception'': 1, ''canbus,drivers'': 1}'
$ws.Range("H4").Value = 2190

$ws.Range("G5").Value = '{''perception,prediction'': 1, ''canbus,perception'': 2, ''drivers,localization,perception'': 1, ''canbus,common,drivers,localization,perception,planning,prediction'': 1, ''common,control,localization,perception,planning,prediction'': 2, ''canbus,common,control,dreamview,localization,map,perception,planning,prediction'': 3, ''common,control,dreamview,localization,perception,planning,prediction'': 4, ''common,control,dreamview,localization,map,perception,planning,prediction'': 2, ''drivers,localization,map,perception'': 6, ''control,localization,perception'': 1, ''planning,prediction'': 8, ''localization,planning,prediction'': 1, ''map,planning,prediction'': 1, ''perception,planning,prediction'': 2, ''control,localization'': 1, ''control,localization,planning,prediction'': 1, ''common,localization'': 2, ''map,perception'': 3, ''common,localization,perception,planning,prediction'': 2, ''drivers,localization'': 7, ''perception,planning'': 4, ''drivers,perception'': 2, ''dreamview,map'': 4, ''dreamview,map,planning,prediction'': 2, ''localization,perception,planning,prediction'': 1, ''control,localization,perception,planning,prediction'': 1, ''control,localization,perception,planning'': 1, ''perception,planning,prediction,routing'': 1, ''common,map'': 3, ''localization,perception'': 2, ''map,prediction'': 1, ''control,map'': 1, ''map,planning'': 1, ''control,planning'': 3, ''canbus,common,control,dreamview,drivers,localization,map,perception,planning,prediction,routing'': 1, ''control,localization,planning'': 1, ''localization,planning'': 3, ''control,prediction'': 1, ''canbus,drivers'': 1, ''dreamview,localization'': 1}'
$ws.Range("H5").Value = 2239

$ws.Range("G6").Value = '{''planning,prediction'': 6, ''map,planning'': 2, ''dreamview,map,planning,prediction'': 2, ''dreamview,map'': 5, ''control,localization'': 1, ''drivers,planning'': 1, ''control,planning'': 3, ''control,localization,planning,prediction'': 2, ''localization,planning,prediction'': 1, ''perception,planning,prediction'': 3, ''control,prediction'': 2, ''perception,prediction'': 1, ''canbus,common,control,dreamview,drivers,localization,map,perception,planning,prediction,routing'': 1, ''canbus,common,control,dreamview,localization,map,perception,planning,prediction'': 1, ''canbus,common,control,dreamview,drivers,localization,map,perception,planning,prediction'': 1, ''control,localization,planning'': 1, ''drivers,localization'': 4, ''drivers,localization,map'': 5, ''planning,prediction,routing'': 1, ''common,drivers'': 2, ''common,perception'': 1, ''canbus,common,control,dreamview,localization,perception,planning,prediction'': 5, ''common,control,dreamview,localization,perception,planning,prediction'': 7, ''common,dreamview,localization,perception,planning,prediction'': 1, ''common,localization,perception,planning,prediction'': 2, ''control,drivers,map'': 1, ''control,drivers'': 1, ''perception,planning'': 3, ''localization,perception,planning,prediction'': 1, ''localization,planning'': 3, ''common,map'': 4, ''control,map,prediction'': 1, ''common,localization'': 1, ''control,map'': 1, ''dreamview,localization'': 1, ''localization,map'': 1}'
$ws.Range("H6").Value = 1914

$ws.Range("G7").Value = '{''perception,v2x'': 2, ''drivers,localization'': 2, ''dreamview,localization'': 1, ''planning,prediction'': 12, ''localization,perception,planning,prediction'': 1, ''perception,planning,prediction'': 1, ''perception,prediction'': 5, ''canbus,common,control,dreamview,drivers,localization,map,perception,planning,prediction'': 1, ''drivers,localization,map,perception,planning'': 2, ''drivers,localization,map,perception'': 1, ''control,planning'': 9, ''perception,planning'': 16, ''control,localization,perception,planning'': 1, ''control,perception,planning'': 1, ''monitor,perception'': 1, ''canbus,common,control,drivers,localization,map,perception,planning,prediction,v2x'': 2, ''map,planning'': 1, ''common,map'': 2, ''canbus,common,control,localization,perception,planning,prediction,v2x'': 2, ''common,perception'': 5, ''canbus,common,control,localization,perception,planning,prediction'': 2, ''canbus,common,control,dreamview,localization,perception,planning,prediction,v2x'': 1, ''common,control,localization,perception,planning,prediction'': 1, ''common,control,dreamview,localization,perception,planning,prediction'': 2, ''common,localization,perception,planning,prediction'': 1, ''common,localization,planning,prediction'': 1, ''map,v2x'': 1, ''drivers,planning'': 1, ''drivers,perception'': 2, ''localization,planning,prediction'': 1, ''drivers,map'': 2, ''dreamview,map,perception,planning,prediction'': 1, ''dreamview,map,perception,prediction'': 1, ''dreamview,map,prediction'': 1, ''dreamview,map'': 1, ''control,localization,map'': 2, ''localization,planning'': 2, ''localization,map'': 2, ''common,planning'': 2, ''control,map'': 1, ''control,perception'': 4, ''control,localization'': 1, ''map,perception'': 3, ''drivers,localization,routing'': 1, ''canbus,control'': 2, ''drivers,prediction'': 1, ''localization,perception'': 3, ''localization,tools'': 2, ''control,drivers'': 1, ''map,tools'': 1, ''dreamview,perception'': 2}'
$ws.Range("H7").Value = 1642

$ws.Range("G8").Value = '{''canbus,common,control,data,dreamview,drivers,localization,map,perception,planning,prediction'': 1, ''common,prediction'': 4, ''control,perception'': 4, ''control,perception,planning'': 4, ''perception,planning'': 10, ''guardian,monitor'': 1, ''dreamview,prediction'': 1, ''common,control,dreamview,localization,perception,planning,prediction'': 5, ''canbus,common,control,dreamview,localization,perception,planning,prediction,v2x'': 2, ''planning,prediction'': 9, ''localization,planning'': 1, ''drivers,planning'': 3, ''control,localization,perception,planning,prediction'': 1, ''localization,tools'': 2, ''map,tools'': 1, ''localization,perception'': 6, ''control,drivers,perception,planning'': 1, ''control,planning'': 6, ''canbus,common,control,localization,perception,planning'': 1, ''common,control,localization,perception,planning,prediction'': 2, ''common,localization,perception,planning,prediction'': 1, ''common,localization,planning,prediction'': 1, ''common,planning,prediction'': 1, ''common,perception'': 1, ''map,perception'': 7, ''canbus,drivers'': 3, ''canbus,common,control,drivers,localization,map,perception,planning,v2x'': 1, ''map,planning'': 3, ''canbus,common,control,drivers,localization,perception,planning,v2x'': 1, ''canbus,common,control,localization,perception,planning,v2x'': 3, ''localization,map'': 1, ''dreamview,map,perception,planning'': 1, ''dreamview,map'': 3, ''drivers,perception'': 1, ''map,v2x'': 1, ''dreamview,map,perception'': 1, ''localization,map,perception'': 2, ''drivers,localization,map,perception'': 1, ''drivers,localization'': 2, ''dreamview,planning'': 1, ''dreamview,perception'': 1, ''drivers,map'': 1, ''map,prediction'': 1, ''common,planning'': 1}'
$ws.Range("H8").Value = 1696
